# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn and de-de sheets to
# reflect that the b.md handoff is now ready (new handoff file/date),
# replacing the previous "Handed back: in sync with en-US" status.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 corresponds to b.md ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-26-19 18:26:48"

# --- zh-cn sheet: row 3 corresponds to b.md ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 18:26:45"

# --- de-de sheet: row 3 corresponds to b.md ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 18:26:48"
